# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# for the last row (b7277af2-cbe8-4e42-9c08-0f7360d1714a.md) on each sheet to reflect
# a fresh handoff report generation.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G5").Value = "2016-11-07 06:51:08"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H5").Value = "2016-11-07 06:50:54"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H5").Value = "2016-11-07 06:51:08"
